$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new test case row (row 8)
$ws.Range("A8").Value = "longest word ends with apostrophe"
$ws.Range("B8").Value = "The big words' end."
$ws.Range("C8").Value = 6
$ws.Range("D8").Value = "words'"

# Select the new last cell, matching the saved view state
$ws.Range("D8").Select()
